# Update "想去人数" (want-to-go count) values in both the "展览" sheet
# and the "全部类型" sheet, which both contain the same underlying events.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 522
$ws1.Range("F4").Value = 188
$ws1.Range("F7").Value = 94
$ws1.Range("F8").Value = 107
$ws1.Range("F9").Value = 39
$ws1.Range("F10").Value = 6563
$ws1.Range("F11").Value = 223
$ws1.Range("F12").Value = 356
$ws1.Range("F13").Value = 2760
$ws1.Range("F14").Value = 171
$ws1.Range("F15").Value = 289
$ws1.Range("F17").Value = 516

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 522
$ws4.Range("F6").Value = 188
$ws4.Range("F9").Value = 94
$ws4.Range("F10").Value = 107
$ws4.Range("F11").Value = 39
$ws4.Range("F13").Value = 6563
$ws4.Range("F15").Value = 223
$ws4.Range("F16").Value = 356
$ws4.Range("F17").Value = 2760
$ws4.Range("F18").Value = 171
$ws4.Range("F19").Value = 289
$ws4.Range("F21").Value = 516
